$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new data row (row 13) that was previously blank.
$ws.Range("B13").Value = 3077271
$ws.Range("C13").Value = 70

# Re-enter/fill the "Tx cost (Eth)" formula down column F so it becomes a
# shared formula group spanning F4:F17 (mirrors re-filling the column after
# adding the new row of data).
$ws.Range("F4:F17").Formula = '=IF(D4<>"",D4*0.000000001,"")'

# Move the active selection to C14, matching where the user continued entry.
$ws.Range("C14").Select() | Out-Null
